# BOM-REV3.xlsx update: add missing components (C29, D8, D9, R27) and
# merge U7-U8 into U7-U9 (DGTL ISO part now covers 3 designators).
#
# The sheet is sorted by "Customer Reference" (column F), so the new
# components are inserted at the correct alphabetical/sorted position,
# matching how Excel would look after the author re-sorted the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert blank rows at their original (pre-edit) positions, working from
# the bottom of the sheet upward so that earlier inserts don't shift the
# row numbers used by later inserts.

# 1) New row for R27 (120K 1% resistor) - was inserted just above old row 45 (R33)
$ws.Rows.Item(45).Insert()

# 2) New rows for D9 and D8 (diodes) - inserted just above old row 20 (D5 D12)
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# 3) New row for C29 (0.1uF 0805 cap) - inserted just above old row 9 (C17 C24)
$ws.Rows.Item(9).Insert()

# Now fill in the data for the newly inserted rows (final row numbers,
# after all four inserts above have shifted everything below them down).

# Row 9: C29
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = "399-1169-1-ND"
$ws.Cells.Item(9, 4).Value = "C0805C104M5RACTU"
$ws.Cells.Item(9, 5).Value = "CAP CER 0.1UF 50V X7R 0805"
$ws.Cells.Item(9, 6).Value = "C29"
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0.11
$ws.Cells.Item(9, 10).Value = 0.11

# Row 21: D9
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(21, 3).Value = "MBR0520LCT-ND"
$ws.Cells.Item(21, 4).Value = "MBR0520L"
$ws.Cells.Item(21, 5).Value = "DIODE SCHOTTKY 20V 500MA SOD123"
$ws.Cells.Item(21, 6).Value = "D9"
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0.34
$ws.Cells.Item(21, 10).Value = 0.34

# Row 22: D8
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = "SMA6J18AHR3GCT-ND"
$ws.Cells.Item(22, 4).Value = "SMA6J18AHR3G"
$ws.Cells.Item(22, 5).Value = "TVS DIODE 18V 28.3V DO214AC"
$ws.Cells.Item(22, 6).Value = "D8"
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0.47
$ws.Cells.Item(22, 10).Value = 0.47

# Row 48: R27
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = 1
$ws.Cells.Item(48, 3).Value = "RMCF0402FT120KCT-ND"
$ws.Cells.Item(48, 4).Value = "RMCF0402FT120K"
$ws.Cells.Item(48, 5).Value = "RES 120K OHM 1% 1/16W 0402"
$ws.Cells.Item(48, 6).Value = "R27"
$ws.Cells.Item(48, 7).Value = 1
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0.1
$ws.Cells.Item(48, 10).Value = 0.1

# Row 57 (was the "U7-U8" row, now shifted down by the 4 inserts above):
# the part now covers U7, U8 and U9, so only the reference designator
# label changes.
$ws.Cells.Item(57, 6).Value = "U7-U9"

# Update the index column (A) for every data row so it again reads
# 1, 2, 3, ... (row - 1), since the inserts only auto-filled the rows
# that moved, not the brand new ones outside that range.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Refresh the totals row so it sums the full, now-larger, data range.
$totalRow = $lastRow + 1
$ws.Cells.Item($totalRow, 7).Formula = "=SUM(G2:G$lastRow)"
$ws.Cells.Item($totalRow, 10).Formula = "=SUM(J2:J$lastRow)"

# Restore the view state (selection) roughly where the author left it.
$ws.Range("F55").Select()
